$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Freelancing/5000 -> Salary/200000, date updated
$ws.Range("A2").Value = "Salary"
$ws.Range("B2").Value = 200000
$ws.Range("C2").Value = 45809.22928240741

# Row 3: Freelancing/5000 -> Ola/3000, date updated
$ws.Range("A3").Value = "Ola"
$ws.Range("B3").Value = 3000
$ws.Range("C3").Value = 45770.22928240741

# Row 4: new row - Food Business/1000/date.
# Insert the row first so the new C4 cell inherits the C3 date style
# (same numFmt as the other date cells) rather than getting a fresh style.
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Food Business"
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 45764.22928240741
